$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new row 40 with the latest EUR->ARS quote.
# Use explicit text-formatted cells so values stay as strings (matching
# the existing rows), not auto-converted to dates/numbers by Excel.
$row = 40

$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2025-09-25"

$cellB = $ws.Cells.Item($row, 2)
$cellB.Value = "21:20:03"

$cellC = $ws.Cells.Item($row, 3)
$cellC.Value = "1.00 EUR = 1,629.0690"
